# Fruta / hortaliza, semanal
# Re-order the weekly price records (rows 2-8) to reflect the latest
# weekly data pull. Columns A,B,C,E,F,G,H,I,J are identical across all
# rows, so only D (Fecha), K (Variedad), L (Calidad), M (Volumen),
# N (Precio minimo), O (Precio maximo), P (Precio promedio ponderado),
# Q (Unidad de comercializacion), R (Origen), S (Precio $/Kg) and
# T (Kg / unidad) need to be rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{
    2 = @{ D = 44229; K = "Santina"; L = "Primera"; M = 250; N = 6500;  O = 7000;  P = 6750;  Q = "`$/bandeja 5 kilos";  R = "Provincia de Curicó";  S = 1350; T = 5 }
    3 = @{ D = 44537; K = "Brooks";  L = "Primera"; M = 200; N = 29000; O = 30000; P = 29500; Q = "`$/caja 20 kilos";     R = "Región de O'Higgins"; S = 1475; T = 20 }
    4 = @{ D = 44532; K = "Brooks";  L = "Primera"; M = 400; N = 27000; O = 28000; P = 27500; Q = "`$/bandeja 12 kilos"; R = "Región de O'Higgins"; S = 2292; T = 12 }
    5 = @{ D = 44210; K = "Rainier"; L = "Segunda"; M = 250; N = 21000; O = 22000; P = 21500; Q = "`$/caja 18 kilos";     R = "Región de O'Higgins"; S = 1194; T = 18 }
    6 = @{ D = 44161; K = "Bing";    L = "Primera"; M = 160; N = 39000; O = 40000; P = 39500; Q = "`$/caja 20 kilos";     R = "Provincia de Curicó";  S = 1975; T = 20 }
    7 = @{ D = 44208; K = "Lapins";  L = "Segunda"; M = 200; N = 10500; O = 11000; P = 10750; Q = "`$/bandeja 12 kilos"; R = "Provincia de Curicó";  S = 896;  T = 12 }
    8 = @{ D = 44175; K = "Rainier"; L = "Segunda"; M = 270; N = 25000; O = 26000; P = 25500; Q = "`$/caja 18 kilos";     R = "Región de O'Higgins"; S = 1417; T = 18 }
}

foreach ($r in $rows.Keys) {
    $data = $rows[$r]
    $ws.Range("D$r").Value = $data.D
    $ws.Range("K$r").Value = $data.K
    $ws.Range("L$r").Value = $data.L
    $ws.Range("M$r").Value = $data.M
    $ws.Range("N$r").Value = $data.N
    $ws.Range("O$r").Value = $data.O
    $ws.Range("P$r").Value = $data.P
    $ws.Range("Q$r").Value = $data.Q
    $ws.Range("R$r").Value = $data.R
    $ws.Range("S$r").Value = $data.S
    $ws.Range("T$r").Value = $data.T
}
